$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# Header text updates (rich-text shared strings resolved to their full concatenated text)
$ws.Range("A8").Value = "Volume 32   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/21/2025  Through  7/27/2025"

# Crime statistics table updates (rows 14-31, 33)
$ws.Range("D14").Value = 3
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = -57.142857142857
$ws.Range("J14").Value = 46
$ws.Range("K14").Value = -39.130434782608
$ws.Range("M14").Value = -65
$ws.Range("N14").Value = -89.433962264150
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = -75
$ws.Range("F15").Value = 21
$ws.Range("G15").Value = 23
$ws.Range("H15").Value = -8.695652173913
$ws.Range("I15").Value = 148
$ws.Range("J15").Value = 148
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 10.447761194029
$ws.Range("M15").Value = 25.423728813559
$ws.Range("N15").Value = -58.659217877095
$ws.Range("C16").Value = 41
$ws.Range("E16").Value = -14.583333333333
$ws.Range("F16").Value = 184
$ws.Range("G16").Value = 179
$ws.Range("H16").Value = 2.793296089385
$ws.Range("I16").Value = 1207
$ws.Range("J16").Value = 1376
$ws.Range("K16").Value = -12.281976744186
$ws.Range("L16").Value = -11.897810218978
$ws.Range("M16").Value = -37.168141592920
$ws.Range("N16").Value = -86.851851851851
$ws.Range("C17").Value = 71
$ws.Range("D17").Value = 108
$ws.Range("E17").Value = -34.259259259259
$ws.Range("F17").Value = 329
$ws.Range("G17").Value = 399
$ws.Range("H17").Value = -17.543859649122
$ws.Range("I17").Value = 2488
$ws.Range("J17").Value = 2556
$ws.Range("K17").Value = -2.660406885759
$ws.Range("L17").Value = 0.120724346076
$ws.Range("M17").Value = 30.057501306847
$ws.Range("N17").Value = -48.637489677952
$ws.Range("C18").Value = 31
$ws.Range("D18").Value = 42
$ws.Range("E18").Value = -26.190476190476
$ws.Range("F18").Value = 122
$ws.Range("G18").Value = 139
$ws.Range("H18").Value = -12.230215827338
$ws.Range("I18").Value = 903
$ws.Range("J18").Value = 1083
$ws.Range("K18").Value = -16.620498614958
$ws.Range("L18").Value = -22.820512820512
$ws.Range("M18").Value = -45.569620253164
$ws.Range("N18").Value = -86.661742983751
$ws.Range("C19").Value = 107
$ws.Range("D19").Value = 130
$ws.Range("E19").Value = -17.692307692307
$ws.Range("F19").Value = 458
$ws.Range("G19").Value = 478
$ws.Range("H19").Value = -4.184100418410
$ws.Range("I19").Value = 2987
$ws.Range("J19").Value = 3083
$ws.Range("K19").Value = -3.113850145961
$ws.Range("L19").Value = -10.219416892095
$ws.Range("M19").Value = 31.934628975265
$ws.Range("N19").Value = -19.919571045576
$ws.Range("C20").Value = 50
$ws.Range("D20").Value = 43
$ws.Range("E20").Value = 16.279069767441
$ws.Range("F20").Value = 183
$ws.Range("G20").Value = 146
$ws.Range("H20").Value = 25.342465753424
$ws.Range("I20").Value = 892
$ws.Range("J20").Value = 959
$ws.Range("K20").Value = -6.986444212721
$ws.Range("L20").Value = -11.770524233432
$ws.Range("M20").Value = 9.717097170971
$ws.Range("N20").Value = -83.376816995900
$ws.Range("C21").Value = 302
$ws.Range("D21").Value = 382
$ws.Range("E21").Value = -20.942408376963
$ws.Range("F21").Value = 1300
$ws.Range("G21").Value = 1371
$ws.Range("H21").Value = -5.178701677607
$ws.Range("I21").Value = 8653
$ws.Range("J21").Value = 9251
$ws.Range("K21").Value = -6.464166036104
$ws.Range("L21").Value = -9.250131095962
$ws.Range("M21").Value = -1.311587591240
$ws.Range("N21").Value = -71.641595385573
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = -80
$ws.Range("F22").Value = 21
$ws.Range("G22").Value = 27
$ws.Range("H22").Value = -22.222222222222
$ws.Range("I22").Value = 175
$ws.Range("J22").Value = 170
$ws.Range("K22").Value = 2.941176470588
$ws.Range("L22").Value = -2.234636871508
$ws.Range("M22").Value = -23.913043478260
$ws.Range("C23").Value = 19
$ws.Range("D23").Value = 33
$ws.Range("E23").Value = -42.424242424242
$ws.Range("F23").Value = 88
$ws.Range("G23").Value = 122
$ws.Range("H23").Value = -27.868852459016
$ws.Range("I23").Value = 804
$ws.Range("J23").Value = 864
$ws.Range("K23").Value = -6.944444444444
$ws.Range("L23").Value = -14.102564102564
$ws.Range("M23").Value = 25.821596244131
$ws.Range("C24").Value = 285
$ws.Range("D24").Value = 263
$ws.Range("E24").Value = 8.365019011406
$ws.Range("F24").Value = 1099
$ws.Range("G24").Value = 978
$ws.Range("H24").Value = 12.372188139059
$ws.Range("I24").Value = 7140
$ws.Range("J24").Value = 6891
$ws.Range("K24").Value = 3.613408794079
$ws.Range("L24").Value = 0.918727915194
$ws.Range("M24").Value = 24.803356056633
$ws.Range("C25").Value = 101
$ws.Range("D25").Value = 112
$ws.Range("E25").Value = -9.821428571428
$ws.Range("F25").Value = 395
$ws.Range("G25").Value = 465
$ws.Range("H25").Value = -15.053763440860
$ws.Range("I25").Value = 2739
$ws.Range("J25").Value = 3062
$ws.Range("K25").Value = -10.548661005878
$ws.Range("L25").Value = -3.590285110876
$ws.Range("C26").Value = 111
$ws.Range("D26").Value = 136
$ws.Range("E26").Value = -18.382352941176
$ws.Range("F26").Value = 532
$ws.Range("G26").Value = 492
$ws.Range("H26").Value = 8.130081300813
$ws.Range("I26").Value = 3531
$ws.Range("J26").Value = 3611
$ws.Range("K26").Value = -2.215452783162
$ws.Range("L26").Value = 1.116838487972
$ws.Range("M26").Value = -23.122142390594
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = -70
$ws.Range("F27").Value = 24
$ws.Range("G27").Value = 28
$ws.Range("H27").Value = -14.285714285714
$ws.Range("I27").Value = 184
$ws.Range("J27").Value = 214
$ws.Range("K27").Value = -14.018691588785
$ws.Range("L27").Value = -10.243902439024
$ws.Range("C28").Value = 13
$ws.Range("D28").Value = 13
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 59
$ws.Range("G28").Value = 56
$ws.Range("H28").Value = 5.357142857142
$ws.Range("I28").Value = 403
$ws.Range("J28").Value = 361
$ws.Range("K28").Value = 11.634349030470
$ws.Range("L28").Value = 16.138328530259
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 7
$ws.Range("E29").Value = -28.571428571428
$ws.Range("G29").Value = 27
$ws.Range("H29").Value = -29.629629629629
$ws.Range("I29").Value = 113
$ws.Range("J29").Value = 144
$ws.Range("K29").Value = -21.527777777777
$ws.Range("L29").Value = -12.403100775193
$ws.Range("M29").Value = -61.168384879725
$ws.Range("N29").Value = -89.642529789184
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = 4
$ws.Range("F30").Value = 17
$ws.Range("G30").Value = 20
$ws.Range("H30").Value = -15
$ws.Range("I30").Value = 94
$ws.Range("J30").Value = 122
$ws.Range("K30").Value = -22.950819672131
$ws.Range("L30").Value = -16.071428571428
$ws.Range("M30").Value = -59.482758620689
$ws.Range("N30").Value = -90.417940876656
$ws.Range("D31").Value = 4
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = -50
$ws.Range("I31").Value = 50
$ws.Range("J31").Value = 44
$ws.Range("K31").Value = 13.636363636363
$ws.Range("L31").Value = 13.636363636363
$ws.Range("F33").Value = 3
